# add - management command - generate excel dummy transactions
#
# Populates the sale-invoice-detail seed sheet with a full set of dummy
# transaction rows (32 sale invoices, each with 1-3 product-variation
# lines and a random quantity).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("SI-1","PV-001",13),
    @("SI-2","PV-001",17),
    @("SI-2","PV-001",11),
    @("SI-2","PV-002",13),
    @("SI-3","PV-001",1),
    @("SI-3","PV-003",6),
    @("SI-3","PV-003",11),
    @("SI-4","PV-002",11),
    @("SI-5","PV-002",11),
    @("SI-6","PV-003",9),
    @("SI-6","PV-003",18),
    @("SI-6","PV-001",12),
    @("SI-7","PV-001",13),
    @("SI-7","PV-001",3),
    @("SI-8","PV-001",16),
    @("SI-9","PV-002",15),
    @("SI-9","PV-002",15),
    @("SI-9","PV-001",20),
    @("SI-10","PV-001",15),
    @("SI-10","PV-003",5),
    @("SI-11","PV-003",9),
    @("SI-11","PV-003",9),
    @("SI-12","PV-003",17),
    @("SI-12","PV-002",16),
    @("SI-12","PV-003",16),
    @("SI-13","PV-002",14),
    @("SI-13","PV-002",14),
    @("SI-13","PV-001",3),
    @("SI-14","PV-001",4),
    @("SI-15","PV-003",1),
    @("SI-15","PV-003",4),
    @("SI-16","PV-002",19),
    @("SI-16","PV-001",11),
    @("SI-17","PV-002",8),
    @("SI-17","PV-003",19),
    @("SI-17","PV-001",17),
    @("SI-18","PV-001",12),
    @("SI-18","PV-003",8),
    @("SI-19","PV-003",15),
    @("SI-19","PV-001",2),
    @("SI-20","PV-002",11),
    @("SI-21","PV-001",13),
    @("SI-21","PV-002",12),
    @("SI-21","PV-001",19),
    @("SI-22","PV-003",18),
    @("SI-22","PV-001",4),
    @("SI-22","PV-001",12),
    @("SI-23","PV-001",5),
    @("SI-24","PV-003",15),
    @("SI-25","PV-003",3),
    @("SI-25","PV-001",20),
    @("SI-26","PV-001",16),
    @("SI-27","PV-003",15),
    @("SI-27","PV-003",8),
    @("SI-27","PV-002",15),
    @("SI-28","PV-001",14),
    @("SI-28","PV-002",1),
    @("SI-29","PV-002",13),
    @("SI-29","PV-001",8),
    @("SI-30","PV-001",15),
    @("SI-30","PV-001",16),
    @("SI-31","PV-001",9),
    @("SI-32","PV-003",12)
)

$row = 2
foreach ($line in $data) {
    $ws.Cells.Item($row, 1).Value = $line[0]
    $ws.Cells.Item($row, 2).Value = $line[1]
    $ws.Cells.Item($row, 3).Value = $line[2]
    $row++
}

$lastRow = $row - 1
$ws.Range("A2:C$lastRow").Select()
